# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain the same set of events (mirrored data), so the same
# row -> new-value updates need to be applied on each sheet, using the
# row numbers that are specific to each sheet's layout.

$wb = $excel.ActiveWorkbook

# Row numbers are for the "展览" sheet (sheetId 1)
$updatesSheet1 = @{
    2  = 46
    3  = 21320
    4  = 817
    8  = 7949
    9  = 556
    11 = 764
    15 = 172
    16 = 33
    18 = 228
    20 = 538
    24 = 82
    26 = 354
    30 = 226
    34 = 140
    35 = 5068
    38 = 47
    40 = 13144
    41 = 1370
    42 = 139
    43 = 54
    46 = 439
    47 = 4065
}

# Row numbers are for the "全部类型" sheet (sheetId 4)
$updatesSheet4 = @{
    2  = 46
    3  = 21320
    4  = 817
    7  = 7949
    8  = 556
    10 = 764
    14 = 172
    15 = 33
    16 = 228
    18 = 538
    22 = 82
    24 = 354
    28 = 226
    33 = 140
    35 = 5068
    38 = 47
    40 = 13144
    41 = 1370
    42 = 139
    43 = 54
    46 = 439
    47 = 4065
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updatesSheet1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesSheet1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesSheet4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesSheet4[$row]
}
